$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new sample row was recorded and prepended right after the header
# (row 2), pushing the existing gyroscope readings down by one row and
# dropping the oldest reading that used to sit at the bottom of the
# sheet (old row 21) off the used range.
#
# Shift every existing data row (2..21) down into (3..22) by copying
# values upward-to-downward starting from the bottom so we never
# clobber a row before it has been copied.
for ($r = 21; $r -ge 2; $r--) {
    for ($c = 1; $c -le 3; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($r + 1, $c)
        $dstCell.Value2 = $srcCell.Value()
    }
}

# Write the new first observation into row 2.
$ws.Cells.Item(2, 1).Value2 = 0.08162501163598926
$ws.Cells.Item(2, 2).Value2 = -0.581960884536185
$ws.Cells.Item(2, 3).Value2 = 0.1855354215495473

# The shift above duplicated the former last row (old row 21, now at
# row 22) — remove that trailing duplicate so the sheet ends at row 21
# (A1:C21), matching the new, smaller used range.
$ws.Rows.Item(22).Delete()
